$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 5716.75
$ws.Range("I51").Value = 2125.25
$ws.Range("J51").Value = 7512.5
$ws.Range("K51").Value = 2125.25
$ws.Range("L51").Value = 7512.5
$ws.Range("M51").Value = -1641.25
$ws.Range("N51").Value = -8480.5
$ws.Range("H69").Value = 15537.625
$ws.Range("I69").Value = 10716.143
$ws.Range("J69").Value = 17522.941
$ws.Range("K69").Value = 32148.429
$ws.Range("L69").Value = 52568.823
$ws.Range("M69").Value = -31274.429
$ws.Range("N69").Value = -54316.823
$ws.Range("H70").Value = 1363.375
$ws.Range("I70").Value = 902
$ws.Range("J70").Value = 1517.1666
$ws.Range("K70").Value = 2706
$ws.Range("L70").Value = 4551.4998
$ws.Range("M70").Value = -2436
$ws.Range("N70").Value = -5091.4998
$ws.Range("H72").Value = 15537.625
$ws.Range("I72").Value = 10716.143
$ws.Range("J72").Value = 17522.941
$ws.Range("K72").Value = 96445.287
$ws.Range("L72").Value = 157706.469
$ws.Range("M72").Value = -92077.287
$ws.Range("N72").Value = -166442.469
$ws.Range("H73").Value = 1363.375
$ws.Range("I73").Value = 902
$ws.Range("J73").Value = 1517.1666
$ws.Range("K73").Value = 2706
$ws.Range("L73").Value = 4551.4998
$ws.Range("M73").Value = -1770
$ws.Range("N73").Value = -6423.4998
$ws.Range("H80").Value = 2465.1765
$ws.Range("I80").Value = 1611.3334
$ws.Range("K80").Value = 4834.0002
$ws.Range("M80").Value = -3836.0002
$ws.Range("H82").Value = 1147.1666
$ws.Range("I82").Value = 1147.1666
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3441.4998
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -3035.4998
$ws.Range("N82").ClearContents()
$ws.Range("H83").Value = 2465.1765
$ws.Range("I83").Value = 1611.3334
$ws.Range("K83").Value = 14502.0006
$ws.Range("M83").Value = -9510.000599999999
$ws.Range("H85").Value = 1147.1666
$ws.Range("I85").Value = 1147.1666
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3441.4998
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2037.4998
$ws.Range("N85").ClearContents()
$ws.Range("H88").Value = 9699.286
$ws.Range("J88").Value = 9578.4
$ws.Range("L88").Value = 9578.4
$ws.Range("N88").Value = -10390.4
$ws.Range("H91").Value = 9699.286
$ws.Range("J91").Value = 9578.4
$ws.Range("L91").Value = 9578.4
$ws.Range("N91").Value = -12386.4
$ws.Range("H92").Value = 708.1429000000001
$ws.Range("I92").Value = 548.3889
$ws.Range("J92").Value = 1666.6666
$ws.Range("K92").Value = 548.3889
$ws.Range("L92").Value = 1666.6666
$ws.Range("M92").Value = 699.6111
$ws.Range("N92").Value = -4162.6666
$ws.Range("H107").Value = 574.3125
$ws.Range("I107").Value = 585.5333000000001
$ws.Range("J107").Value = 406
$ws.Range("K107").Value = 585.5333000000001
$ws.Range("L107").Value = 406
$ws.Range("M107").Value = 1334.4667
$ws.Range("N107").Value = -4246
$ws.Range("H111").Value = 5782.636
$ws.Range("I111").Value = 4532.7144
$ws.Range("J111").Value = 7970
$ws.Range("K111").Value = 13598.1432
$ws.Range("L111").Value = 23910
$ws.Range("M111").Value = -10531.1432
$ws.Range("N111").Value = -30044

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 1889.6923
$ws.Range("I88").Value = 1295.1428
$ws.Range("K88").Value = 1295.1428
$ws.Range("M88").Value = -889.1428000000001
$ws.Range("H91").Value = 1889.6923
$ws.Range("I91").Value = 1295.1428
$ws.Range("K91").Value = 1295.1428
$ws.Range("M91").Value = 108.8571999999999
$ws.Range("H132").Value = 1125.2632
$ws.Range("I132").Value = 902.06665
$ws.Range("J132").Value = 1962.25
$ws.Range("K132").Value = 2706.19995
$ws.Range("L132").Value = 5886.75
$ws.Range("M132").Value = -176.1999500000002
$ws.Range("N132").Value = -10946.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 831.3871
$ws.Range("I94").Value = 659.2727
$ws.Range("J94").Value = 1252.1111
$ws.Range("K94").Value = 659.2727
$ws.Range("L94").Value = 1252.1111
$ws.Range("M94").Value = -208.2727
$ws.Range("N94").Value = -2154.1111

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1281.5862
$ws.Range("I58").Value = 1210.9565
$ws.Range("J58").Value = 1552.3334
$ws.Range("K58").Value = 1210.9565
$ws.Range("L58").Value = 1552.3334
$ws.Range("M58").Value = -1007.9565
$ws.Range("N58").Value = -1958.3334
$ws.Range("H62").Value = 2292.2856
$ws.Range("I62").Value = 2295
$ws.Range("J62").Value = 2285.5
$ws.Range("K62").Value = 2295
$ws.Range("L62").Value = 2285.5
$ws.Range("M62").Value = -1671
$ws.Range("N62").Value = -3533.5
$ws.Range("H65").Value = 2292.2856
$ws.Range("I65").Value = 2295
$ws.Range("J65").Value = 2285.5
$ws.Range("K65").Value = 11475
$ws.Range("L65").Value = 11427.5
$ws.Range("M65").Value = -8355
$ws.Range("N65").Value = -17667.5
$ws.Range("H94").Value = 1238.1
$ws.Range("I94").Value = 1539.5
$ws.Range("J94").Value = 1162.75
$ws.Range("K94").Value = 1539.5
$ws.Range("L94").Value = 1162.75
$ws.Range("M94").Value = -1088.5
$ws.Range("N94").Value = -2064.75
$ws.Range("H136").Value = 1281.5862
$ws.Range("I136").Value = 1210.9565
$ws.Range("J136").Value = 1552.3334
$ws.Range("K136").Value = 3632.8695
$ws.Range("L136").Value = 4657.0002
$ws.Range("M136").Value = -1082.8695
$ws.Range("N136").Value = -9757.0002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 960.25
$ws.Range("I4").Value = 231
$ws.Range("J4").Value = 1203.3334
$ws.Range("K4").Value = 693
$ws.Range("L4").Value = 3610.0002
$ws.Range("M4").Value = -581
$ws.Range("N4").Value = -3834.0002
$ws.Range("H5").Value = 791.25
$ws.Range("I5").Value = 508.3684
$ws.Range("J5").Value = 1388.4445
$ws.Range("K5").Value = 1525.1052
$ws.Range("L5").Value = 4165.333500000001
$ws.Range("M5").Value = -1413.1052
$ws.Range("N5").Value = -4389.333500000001
$ws.Range("H7").Value = 28571532
$ws.Range("I7").Value = 41.666668
$ws.Range("J7").Value = 50000150
$ws.Range("K7").Value = 125.000004
$ws.Range("L7").Value = 150000450
$ws.Range("M7").Value = -13.000004
$ws.Range("N7").Value = -150000674
$ws.Range("H103").Value = 3707.257
$ws.Range("I103").Value = 1027.1666
$ws.Range("J103").Value = 5105.5654
$ws.Range("K103").Value = 3081.4998
$ws.Range("L103").Value = 15316.6962
$ws.Range("M103").Value = -2202.4998
$ws.Range("N103").Value = -17074.6962
$ws.Range("H122").Value = 559.6070999999999
$ws.Range("I122").Value = 533.65216
$ws.Range("J122").Value = 679
$ws.Range("K122").Value = 4802.869439999999
$ws.Range("L122").Value = 6111
$ws.Range("M122").Value = -2352.869439999999
$ws.Range("N122").Value = -11011
$ws.Range("H135").Value = 791.25
$ws.Range("I135").Value = 508.3684
$ws.Range("J135").Value = 1388.4445
$ws.Range("K135").Value = 4575.3156
$ws.Range("L135").Value = 12496.0005
$ws.Range("M135").Value = -2040.3156
$ws.Range("N135").Value = -17566.0005

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2207.1667
$ws.Range("I80").Value = 2109.4443
$ws.Range("J80").Value = 2500.3333
$ws.Range("K80").Value = 2109.4443
$ws.Range("L80").Value = 2500.3333
$ws.Range("M80").Value = -1111.4443
$ws.Range("N80").Value = -4496.3333
$ws.Range("H83").Value = 2207.1667
$ws.Range("I83").Value = 2109.4443
$ws.Range("J83").Value = 2500.3333
$ws.Range("K83").Value = 10547.2215
$ws.Range("L83").Value = 12501.6665
$ws.Range("M83").Value = -5555.2215
$ws.Range("N83").Value = -22485.6665

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 10515.154
$ws.Range("I68").Value = 22100.4
$ws.Range("J68").Value = 3274.375
$ws.Range("K68").Value = 22100.4
$ws.Range("L68").Value = 3274.375
$ws.Range("M68").Value = -21351.4
$ws.Range("N68").Value = -4772.375
$ws.Range("H71").Value = 10515.154
$ws.Range("I71").Value = 22100.4
$ws.Range("J71").Value = 3274.375
$ws.Range("K71").Value = 110502
$ws.Range("L71").Value = 16371.875
$ws.Range("M71").Value = -106758
$ws.Range("N71").Value = -23859.875
$ws.Range("H82").Value = 1721.6666
$ws.Range("I82").Value = 2258.3333
$ws.Range("J82").Value = 1453.3334
$ws.Range("K82").Value = 2258.3333
$ws.Range("L82").Value = 1453.3334
$ws.Range("M82").Value = -1897.3333
$ws.Range("N82").Value = -2175.3334
$ws.Range("H85").Value = 1721.6666
$ws.Range("I85").Value = 2258.3333
$ws.Range("J85").Value = 1453.3334
$ws.Range("K85").Value = 2258.3333
$ws.Range("L85").Value = 1453.3334
$ws.Range("M85").Value = -1010.3333
$ws.Range("N85").Value = -3949.3334
$ws.Range("H132").Value = 4904.7905
$ws.Range("I132").Value = 5601
$ws.Range("J132").Value = 4299.391
$ws.Range("K132").Value = 16803
$ws.Range("L132").Value = 12898.173
$ws.Range("M132").Value = -14273
$ws.Range("N132").Value = -17958.173
$ws.Range("H136").Value = 1448.8909
$ws.Range("I136").Value = 1059.475
$ws.Range("K136").Value = 3178.425
$ws.Range("M136").Value = -628.4249999999997

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2111.125
$ws.Range("I81").Value = 984.1429000000001
$ws.Range("J81").Value = 10000
$ws.Range("K81").Value = 1968.2858
$ws.Range("L81").Value = 20000
$ws.Range("M81").Value = -907.2858000000001
$ws.Range("N81").Value = -22122
$ws.Range("H84").Value = 2111.125
$ws.Range("I84").Value = 984.1429000000001
$ws.Range("J84").Value = 10000
$ws.Range("K84").Value = 9841.429
$ws.Range("L84").Value = 100000
$ws.Range("M84").Value = -4537.429
$ws.Range("N84").Value = -110608
$ws.Range("H132").Value = 2336.6155
$ws.Range("I132").Value = 2462.4595
$ws.Range("J132").Value = 2026.2
$ws.Range("K132").Value = 7387.3785
$ws.Range("L132").Value = 6078.6
$ws.Range("M132").Value = -4857.3785
$ws.Range("N132").Value = -11138.6
